# Insert a new "Availability" column before the existing "EIC-Code" column
# (O), pushing the EIC-Code values/header into column P. This mirrors the
# author's manual edit: a new empty (centered) column O is introduced for
# every data row, the old O header/value moves to P, and the sheet's
# selection / defined name are refreshed to match the new extent.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pp_list_CZ2")

# --- 1. Header row: swap O1 ("EIC-Code") and P1 ("Availability") --------
$oldO1 = $ws.Cells.Item(1, 15).Value()
$oldP1 = $ws.Cells.Item(1, 16).Value()
$ws.Cells.Item(1, 15).Value = $oldP1
$ws.Cells.Item(1, 16).Value = $oldO1

# --- 2. Data rows 2-74: move any existing EIC-Code value from O to P ----
for ($r = 2; $r -le 74; $r++) {
    $oCell = $ws.Cells.Item($r, 15)
    $val = $oCell.Value()
    if ($val -ne $null -and $val -ne "") {
        $ws.Cells.Item($r, 16).Value = $val
        $oCell.Value = $null
    }
}

# --- 3. Every row 1-74 gets the new column O centered (style "3") -------
#     (matches the workbook's existing "empty but centered" cell style,
#      the same one already used on e.g. M2/N1). Apply to the whole
#      column range in one shot to avoid leaving lots of transient
#      single-property style entries behind.
$oRange = $ws.Range("O1:O74")
$oRange.HorizontalAlignment = -4108
$oRange.VerticalAlignment = -4108

# --- 4. Column widths: new col O narrower, col P keeps its old width ----
$ws.Columns.Item(15).ColumnWidth = 10.166666666666666

# --- 5. Selection moves down one row, as recorded in the sheet view -----
$ws.Range("U10").Select() | Out-Null

# --- 6. The hidden _FilterDatabase name now spans through column P ------
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=pp_list_CZ2!`$A`$1:`$P`$74"
    }
}
